$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 267
$ws.Range("A267").Value = 45766
$ws.Range("B267:E267").NumberFormat = "@"
$ws.Range("B267").Value = 'RAVEMANIA EASTER RAVE'
$ws.Range("C267").Value = 'Schlachthof'
$ws.Range("D267").Value = 'Düsseldorf'
$ws.Hyperlinks.Add($ws.Range("E267"), 'https://www.instagram.com/ravemania.official?igsh=aDJsbXBtNHp0aGlu', "", "", 'https://www.instagram.com/ravemania.official?igsh=aDJsbXBtNHp0aGlu')
$ws.Range("E267").Font.Name = "Calibri"
$ws.Range("E267").Font.Size = 11
$ws.Range("E267").Font.Bold = $false
$ws.Range("E267").Font.Underline = $false
$ws.Range("E267").Font.Color = 0
$ws.Range("E267").NumberFormat = "@"

# Row 268
$ws.Range("A268").Value = 45758
$ws.Range("B268:E268").NumberFormat = "@"
$ws.Range("B268").Value = 'CLUB RAVE'
$ws.Range("C268").Value = 'Samy''s Club'
$ws.Range("D268").Value = 'Düsseldorf'
$ws.Hyperlinks.Add($ws.Range("E268"), 'https://www.instagram.com/reel/DHQzw5Noqpi/?igsh=MWMzMzQ2b240c3Bndg==', "", "", 'https://www.instagram.com/reel/DHQzw5Noqpi/?igsh=MWMzMzQ2b240c3Bndg==')
$ws.Range("E268").Font.Name = "Calibri"
$ws.Range("E268").Font.Size = 11
$ws.Range("E268").Font.Bold = $false
$ws.Range("E268").Font.Underline = $false
$ws.Range("E268").Font.Color = 0
$ws.Range("E268").NumberFormat = "@"

# Row 269
$ws.Range("A269").Value = 45745
$ws.Range("B269:E269").NumberFormat = "@"
$ws.Range("B269").Value = 'SAMY''S FAMILY HARD TECHNO'
$ws.Range("C269").Value = 'Samy''s Club'
$ws.Range("D269").Value = 'Düsseldorf'
$ws.Hyperlinks.Add($ws.Range("E269"), 'https://www.instagram.com/p/DHqbJzmsivS/?igsh=MXUydGhuMzRzN25sNg==', "", "", 'https://www.instagram.com/p/DHqbJzmsivS/?igsh=MXUydGhuMzRzN25sNg==')
$ws.Range("E269").Font.Name = "Calibri"
$ws.Range("E269").Font.Size = 11
$ws.Range("E269").Font.Bold = $false
$ws.Range("E269").Font.Underline = $false
$ws.Range("E269").Font.Color = 0
$ws.Range("E269").NumberFormat = "@"

# Row 270
$ws.Range("A270").Value = 45744
$ws.Range("B270:E270").NumberFormat = "@"
$ws.Range("B270").Value = 'RESIDENT NIGHT'
$ws.Range("C270").Value = 'SNRS'
$ws.Range("D270").Value = 'Dortmund'
$ws.Hyperlinks.Add($ws.Range("E270"), 'https://www.instagram.com/reel/DHV6N8ygrzC/?igsh=MTFmbmFxMGFwa2Vyaw==', "", "", 'https://www.instagram.com/reel/DHV6N8ygrzC/?igsh=MTFmbmFxMGFwa2Vyaw==')
$ws.Range("E270").Font.Name = "Calibri"
$ws.Range("E270").Font.Size = 11
$ws.Range("E270").Font.Bold = $false
$ws.Range("E270").Font.Underline = $false
$ws.Range("E270").Font.Color = 0
$ws.Range("E270").NumberFormat = "@"

# Row 271
$ws.Range("A271").Value = 45751
$ws.Range("B271:E271").NumberFormat = "@"
$ws.Range("B271").Value = 'COMMUNITY NIGHT FREE ENTRY'
$ws.Range("C271").Value = 'SNRS'
$ws.Range("D271").Value = 'Dortmund'
$ws.Hyperlinks.Add($ws.Range("E271"), 'https://www.instagram.com/reel/DHtSIQRMmHU/?igsh=ZHh4NmMxNzJ1NXNs', "", "", 'https://www.instagram.com/reel/DHtSIQRMmHU/?igsh=ZHh4NmMxNzJ1NXNs')
$ws.Range("E271").Font.Name = "Calibri"
$ws.Range("E271").Font.Size = 11
$ws.Range("E271").Font.Bold = $false
$ws.Range("E271").Font.Underline = $false
$ws.Range("E271").Font.Color = 0
$ws.Range("E271").NumberFormat = "@"

# Row 272
$ws.Range("A272").Value = 45752
$ws.Range("B272:E272").NumberFormat = "@"
$ws.Range("B272").Value = 'BRUTAL.RADIO'
$ws.Range("C272").Value = 'Die Nacht'
$ws.Range("D272").Value = 'Mönchengladbach'
$ws.Hyperlinks.Add($ws.Range("E272"), 'https://www.instagram.com/reel/DHta3LdtBvi/?igsh=MWR5cDd6a296a3RhZg==', "", "", 'https://www.instagram.com/reel/DHta3LdtBvi/?igsh=MWR5cDd6a296a3RhZg==')
$ws.Range("E272").Font.Name = "Calibri"
$ws.Range("E272").Font.Size = 11
$ws.Range("E272").Font.Bold = $false
$ws.Range("E272").Font.Underline = $false
$ws.Range("E272").Font.Color = 0
$ws.Range("E272").NumberFormat = "@"

# Row 273
$ws.Range("A273").Value = 45766
$ws.Range("B273:E273").NumberFormat = "@"
$ws.Range("B273").Value = 'CLUB TAKEOVER'
$ws.Range("C273").Value = 'Projekt X'
$ws.Range("D273").Value = 'Bochum'
$ws.Hyperlinks.Add($ws.Range("E273"), 'https://www.instagram.com/reel/DHtaJxgoLh0/?igsh=cmFxZ3Nsb3V4aWRr', "", "", 'https://www.instagram.com/reel/DHtaJxgoLh0/?igsh=cmFxZ3Nsb3V4aWRr')
$ws.Range("E273").Font.Name = "Calibri"
$ws.Range("E273").Font.Size = 11
$ws.Range("E273").Font.Bold = $false
$ws.Range("E273").Font.Underline = $false
$ws.Range("E273").Font.Color = 0
$ws.Range("E273").NumberFormat = "@"

# Row 274
$ws.Range("A274").Value = 45759
$ws.Range("B274:E274").NumberFormat = "@"
$ws.Range("B274").Value = 'DIE KINDER DER NACHT'
$ws.Range("C274").Value = 'Stollen134'
$ws.Range("D274").Value = 'Dortmund'
$ws.Hyperlinks.Add($ws.Range("E274"), 'https://www.instagram.com/p/DHQxDzcM5xJ/?igsh=b3llbXFlbTgxYzBm', "", "", 'https://www.instagram.com/p/DHQxDzcM5xJ/?igsh=b3llbXFlbTgxYzBm')
$ws.Range("E274").Font.Name = "Calibri"
$ws.Range("E274").Font.Size = 11
$ws.Range("E274").Font.Bold = $false
$ws.Range("E274").Font.Underline = $false
$ws.Range("E274").Font.Color = 0
$ws.Range("E274").NumberFormat = "@"

# Row 275
$ws.Range("A275").Value = 45752
$ws.Range("B275:E275").NumberFormat = "@"
$ws.Range("B275").Value = 'BASSMANIA MEETS LACUNA'
$ws.Range("C275").Value = 'Favela'
$ws.Range("D275").Value = 'Münster'
$ws.Hyperlinks.Add($ws.Range("E275"), 'https://www.instagram.com/p/DHL5htCtSiD/?igsh=Yml2a3g2YXM4aGRs', "", "", 'https://www.instagram.com/p/DHL5htCtSiD/?igsh=Yml2a3g2YXM4aGRs')
$ws.Range("E275").Font.Name = "Calibri"
$ws.Range("E275").Font.Size = 11
$ws.Range("E275").Font.Bold = $false
$ws.Range("E275").Font.Underline = $false
$ws.Range("E275").Font.Color = 0
$ws.Range("E275").NumberFormat = "@"

# Row 276
$ws.Range("A276").Value = 45759
$ws.Range("B276:E276").NumberFormat = "@"
$ws.Range("B276").Value = 'DAY & NIGHT NEONGREEN MARKET'
$ws.Range("C276").Value = 'Schrotty'
$ws.Range("D276").Value = 'Köln'
$ws.Hyperlinks.Add($ws.Range("E276"), 'https://www.instagram.com/reel/DHrLEMXNTWu/?igsh=MWZqZzV0Mzc0aDE3Yg==', "", "", 'https://www.instagram.com/reel/DHrLEMXNTWu/?igsh=MWZqZzV0Mzc0aDE3Yg==')
$ws.Range("E276").Font.Name = "Calibri"
$ws.Range("E276").Font.Size = 11
$ws.Range("E276").Font.Bold = $false
$ws.Range("E276").Font.Underline = $false
$ws.Range("E276").Font.Color = 0
$ws.Range("E276").NumberFormat = "@"

# Row 277
$ws.Range("A277").Value = 45800
$ws.Range("B277:E277").NumberFormat = "@"
$ws.Range("B277").Value = 'ONE:Z'
$ws.Range("C277").Value = 'Schrotty'
$ws.Range("D277").Value = 'Köln'
$ws.Hyperlinks.Add($ws.Range("E277"), 'https://www.instagram.com/onez.cologne?igsh=YnRnNWl0NDNla2Jv', "", "", 'https://www.instagram.com/onez.cologne?igsh=YnRnNWl0NDNla2Jv')
$ws.Range("E277").Font.Name = "Calibri"
$ws.Range("E277").Font.Size = 11
$ws.Range("E277").Font.Bold = $false
$ws.Range("E277").Font.Underline = $false
$ws.Range("E277").Font.Color = 0
$ws.Range("E277").NumberFormat = "@"

# Row 278
$ws.Range("A278").Value = 45745
$ws.Range("B278:E278").NumberFormat = "@"
$ws.Range("B278").Value = 'HARD RAVE'
$ws.Range("C278").Value = 'Prismatic'
$ws.Range("D278").Value = 'Dortmund'
$ws.Hyperlinks.Add($ws.Range("E278"), 'https://www.instagram.com/p/DHrAKMHMiZd/?igsh=dGh6OTFid2M2NG5v', "", "", 'https://www.instagram.com/p/DHrAKMHMiZd/?igsh=dGh6OTFid2M2NG5v')
$ws.Range("E278").Font.Name = "Calibri"
$ws.Range("E278").Font.Size = 11
$ws.Range("E278").Font.Bold = $false
$ws.Range("E278").Font.Underline = $false
$ws.Range("E278").Font.Color = 0
$ws.Range("E278").NumberFormat = "@"

# Row 279
$ws.Range("A279").Value = 45751
$ws.Range("B279:E279").NumberFormat = "@"
$ws.Range("B279").Value = 'X-BASS'
$ws.Range("C279").Value = 'viersieben'
$ws.Range("D279").Value = 'Duisburg'
$ws.Hyperlinks.Add($ws.Range("E279"), 'https://www.instagram.com/reel/DHlj6UCgVac/?igsh=YXNkNzN6MHR0cGN3', "", "", 'https://www.instagram.com/reel/DHlj6UCgVac/?igsh=YXNkNzN6MHR0cGN3')
$ws.Range("E279").Font.Name = "Calibri"
$ws.Range("E279").Font.Size = 11
$ws.Range("E279").Font.Bold = $false
$ws.Range("E279").Font.Underline = $false
$ws.Range("E279").Font.Color = 0
$ws.Range("E279").NumberFormat = "@"

# Row 280
$ws.Range("A280").Value = 45759
$ws.Range("B280:E280").NumberFormat = "@"
$ws.Range("B280").Value = 'TECHNOBLOCK'
$ws.Range("C280").Value = 'Elektroküche'
$ws.Range("D280").Value = 'Köln'
$ws.Hyperlinks.Add($ws.Range("E280"), 'https://www.instagram.com/reel/DHjJFYusgvO/?igsh=MTRna2xyZGN0bnM0Zw==', "", "", 'https://www.instagram.com/reel/DHjJFYusgvO/?igsh=MTRna2xyZGN0bnM0Zw==')
$ws.Range("E280").Font.Name = "Calibri"
$ws.Range("E280").Font.Size = 11
$ws.Range("E280").Font.Bold = $false
$ws.Range("E280").Font.Underline = $false
$ws.Range("E280").Font.Color = 0
$ws.Range("E280").NumberFormat = "@"

Write-Output "done"
